$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.462
$ws.Range("C6").Value = -12.632
$ws.Range("C7").Value = -13.214
$ws.Range("C8").Value = -12.694
$ws.Range("C16").Value = -12.178
$ws.Range("C20").Value = -12.894
$ws.Range("C21").Value = -13.214
